$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need the column
# pre-formatted as Text so Excel keeps the literal string (matching the
# original workbook, where these price cells are stored as text).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D13", "D14", "D15", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "98.516.24"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.441.51"
$ws.Range("E3").Value = "  +3.35%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "258.94"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "661.76"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("D7").Value = "1.50"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").Value = "0.450"
$ws.Range("E8").Value = "  +7.59%  "
$ws.Range("D9").Value = "1.09"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "3.441.26"
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("D13").Value = "42.68"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").Value = "6.44"
$ws.Range("E14").Value = "  +17.05%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000273"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "98.262.60"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "4.090.92"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("D18").Value = "9.39"
$ws.Range("E18").Value = "  +35.49%  "
$ws.Range("D19").Value = "3.438.05"
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("B20").Value = "Stellar"
$ws.Range("C20").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D20").Value = "0.560"
$ws.Range("E20").Value = "  +31.63%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "17.95"
$ws.Range("E21").Value = "  +7.74%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").Value = "3.54"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  +6.58%  "
$ws.Range("D24").Value = "520.96"
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").Value = "0.0000211"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").Value = "6.45"
$ws.Range("E26").Value = "  +5.03%  "
$ws.Range("D27").Value = "101.89"
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").Value = "13.25"
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("D29").Value = "3.627.47"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").Value = "0.158"
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("D31").Value = "11.89"
$ws.Range("E31").Value = "  +7.42%  "
$ws.Range("D32").Value = "0.203"
$ws.Range("E32").Value = "  +7.74%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "0.592"
$ws.Range("E34").Value = "  +14.48%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "30.26"
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").Value = "  +13.52%  "
$ws.Range("D38").Value = "7.98"
$ws.Range("E38").Value = "  +4.99%  "
$ws.Range("E39").Value = "  +13.38%  "
$ws.Range("D40").Value = "536.64"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("D41").Value = "0.156"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "0.884"
$ws.Range("E43").Value = "  +9.05%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").Value = "9.23"
$ws.Range("E44").Value = "  +20.66%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "24.78"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0440"
$ws.Range("E46").Value = "  +14.49%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "5.84"
$ws.Range("E47").Value = "  +19.35%  "
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").Value = "3.71"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "3.37"
$ws.Range("E49").Value = "  +5.44%  "
$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  +13.88%  "
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").Value = "  +5.58%  "
